$wb = $excel.ActiveWorkbook

# Work on the "addVisitorForExistingIndividual" sheet (2nd sheet).
$ws = $wb.Worksheets.Item("addVisitorForExistingIndividual")

# Populate the new 5x3 block of test data (row by row, left to right so
# shared-string order matches: test, db, q, c, a, w, v, s, e, b, d, r, n, f).
$ws.Range("A1").Value = "lastName"
$ws.Range("B1").Value = "test"
$ws.Range("C1").Value = "db"

$ws.Range("A2").Value = "q"
$ws.Range("B2").Value = "c"
$ws.Range("C2").Value = "a"

$ws.Range("A3").Value = "w"
$ws.Range("B3").Value = "v"
$ws.Range("C3").Value = "s"

$ws.Range("A4").Value = "e"
$ws.Range("B4").Value = "b"
$ws.Range("C4").Value = "d"

$ws.Range("A5").Value = "r"
$ws.Range("B5").Value = "n"
$ws.Range("C5").Value = "f"

# Make this sheet the active tab/selection (matches activeTab=1 and the
# tabSelected flag moving from sheet1 to sheet2 in the workbook).
$ws.Activate()
$ws.Range("A1").Select()
